$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.530.31"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "2.493.63"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'313.85"
$ws.Range("D6").Value = "'93.40"
$ws.Range("E6").Value = "  -1.34%  "
$ws.Range("D7").Value = "'0.547"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'32.70"
$ws.Range("E10").Value = "  -2.81%  "
$ws.Range("D11").Value = "'0.0787"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("E12").Value = "  +1.98%  "
$ws.Range("D13").Value = "2.878.33"
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D15").Value = "'16.20"
$ws.Range("E15").Value = "  +10.44%  "
$ws.Range("D16").Value = "2.500.01"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").Value = "'0.760"
$ws.Range("E17").Value = "  -3.55%  "
$ws.Range("D18").Value = "41.578.77"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("E20").Value = "  +0.89%  "
$ws.Range("D21").Value = "'71.40"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("D23").Value = "'236.42"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("E25").Value = "  -0.55%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "'25.34"
$ws.Range("E27").Value = "  +3.52%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D30").Value = "'36.20"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'158.24"
$ws.Range("E31").Value = "  +2.98%  "
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("D34").Value = "'0.0759"
$ws.Range("E34").Value = "  +0.26%  "
$ws.Range("D35").Value = "'17.99"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("E36").Value = "  -6.17%  "
$ws.Range("D37").Value = "'2.97"
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -1.68%  "
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("E40").Value = "  -0.22%  "
$ws.Range("E41").Value = "  -4.60%  "
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'20.10"
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("D44").Value = "1.969.99"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("E46").Value = "  -2.71%  "
$ws.Range("D47").Value = "'8.90"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "2.727.15"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'96.80"
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("D50").Value = "'68.12"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("D51").Value = "'73.90"
$ws.Range("E51").Value = "  -2.96%  "
